$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1667207.1
$ws.Range("J17").Value = 1667207.1
$ws.Range("L17").Value = 5001621.300000001
$ws.Range("N17").Value = -5001957.300000001
$ws.Range("H46").Value = 23144.389
$ws.Range("J46").Value = 25987.438
$ws.Range("L46").Value = 77962.314
$ws.Range("N46").Value = -78200.314
$ws.Range("H60").Value = 23144.389
$ws.Range("J60").Value = 25987.438
$ws.Range("L60").Value = 77962.314
$ws.Range("N60").Value = -78930.314
$ws.Range("H80").Value = 958.0714
$ws.Range("I80").Value = 166.66667
$ws.Range("J80").Value = 1551.625
$ws.Range("K80").Value = 500.00001
$ws.Range("L80").Value = 4654.875
$ws.Range("M80").Value = 497.99999
$ws.Range("N80").Value = -6650.875
$ws.Range("H83").Value = 958.0714
$ws.Range("I83").Value = 166.66667
$ws.Range("J83").Value = 1551.625
$ws.Range("K83").Value = 1500.00003
$ws.Range("L83").Value = 13964.625
$ws.Range("M83").Value = 3491.99997
$ws.Range("N83").Value = -23948.625
$ws.Range("H87").Value = 24133.107
$ws.Range("J87").Value = 24133.107
$ws.Range("L87").Value = 24133.107
$ws.Range("N87").Value = -26629.107
$ws.Range("H90").Value = 24133.107
$ws.Range("J90").Value = 24133.107
$ws.Range("L90").Value = 72399.321
$ws.Range("N90").Value = -84879.321
$ws.Range("H132").Value = 1000.6129
$ws.Range("I132").Value = 862.3333
$ws.Range("K132").Value = 2586.9999
$ws.Range("M132").Value = -56.9998999999998
$ws.Range("H135").Value = 671.4231
$ws.Range("J135").Value = 1968
$ws.Range("L135").Value = 17712
$ws.Range("N135").Value = -22782
$ws.Range("H138").Value = 1500.0754
$ws.Range("I138").Value = 1271.9286
$ws.Range("J138").Value = 2371.182
$ws.Range("K138").Value = 3815.7858
$ws.Range("L138").Value = 7113.545999999999
$ws.Range("M138").Value = 1324.2142
$ws.Range("N138").Value = -17393.546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2850.875
$ws.Range("J61").Value = 1399.1666
$ws.Range("L61").Value = 1399.1666
$ws.Range("N61").Value = -1823.1666
$ws.Range("H132").Value = 2189.3257
$ws.Range("I132").Value = 1168.3914
$ws.Range("J132").Value = 3363.4
$ws.Range("K132").Value = 3505.1742
$ws.Range("L132").Value = 10090.2
$ws.Range("M132").Value = -975.1741999999999
$ws.Range("N132").Value = -15150.2
$ws.Range("H136").Value = 2850.875
$ws.Range("J136").Value = 1399.1666
$ws.Range("L136").Value = 4197.4998
$ws.Range("N136").Value = -9297.4998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1704.9524
$ws.Range("I134").Value = 969.51514
$ws.Range("J134").Value = 4401.5557
$ws.Range("K134").Value = 2908.54542
$ws.Range("L134").Value = 13204.6671
$ws.Range("M134").Value = -373.5454199999999
$ws.Range("N134").Value = -18274.6671

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2113.3704
$ws.Range("I58").Value = 1780.6471
$ws.Range("J58").Value = 2679
$ws.Range("K58").Value = 1780.6471
$ws.Range("L58").Value = 2679
$ws.Range("M58").Value = -1577.6471
$ws.Range("N58").Value = -3085
$ws.Range("H93").Value = 24571.4
$ws.Range("I93").Value = 5953.5
$ws.Range("J93").Value = 36983.332
$ws.Range("K93").Value = 5953.5
$ws.Range("L93").Value = 36983.332
$ws.Range("M93").Value = -4081.5
$ws.Range("N93").Value = -40727.332
$ws.Range("H132").Value = 3045.85
$ws.Range("I132").Value = 1402
$ws.Range("J132").Value = 4689.7
$ws.Range("K132").Value = 4206
$ws.Range("L132").Value = 14069.1
$ws.Range("M132").Value = -1676
$ws.Range("N132").Value = -19129.1
$ws.Range("H134").Value = 2259.6333
$ws.Range("I134").Value = 2233.2273
$ws.Range("J134").Value = 2332.25
$ws.Range("K134").Value = 6699.6819
$ws.Range("L134").Value = 6996.75
$ws.Range("M134").Value = -4164.6819
$ws.Range("N134").Value = -12066.75
$ws.Range("H136").Value = 2113.3704
$ws.Range("I136").Value = 1780.6471
$ws.Range("J136").Value = 2679
$ws.Range("K136").Value = 5341.9413
$ws.Range("L136").Value = 8037
$ws.Range("M136").Value = -2791.9413
$ws.Range("N136").Value = -13137

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1654.1818
$ws.Range("J5").Value = 3999
$ws.Range("L5").Value = 11997
$ws.Range("N5").Value = -12221
$ws.Range("H132").Value = 1305.1666
$ws.Range("I132").Value = 1305.1666
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11746.4994
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -9216.499400000001
$ws.Range("H135").Value = 1654.1818
$ws.Range("J135").Value = 3999
$ws.Range("L135").Value = 35991
$ws.Range("N135").Value = -41061

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1269.4706
$ws.Range("I2").Value = 1709.1666
$ws.Range("J2").Value = 1029.6364
$ws.Range("K2").Value = 1709.1666
$ws.Range("L2").Value = 1029.6364
$ws.Range("M2").Value = -1596.1666
$ws.Range("N2").Value = -1255.6364

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3220.963
$ws.Range("I136").Value = 3384.3
$ws.Range("J136").Value = 2754.2856
$ws.Range("K136").Value = 10152.9
$ws.Range("L136").Value = 8262.856800000001
$ws.Range("M136").Value = -7602.900000000001
$ws.Range("N136").Value = -13362.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1599.6364
$ws.Range("I132").Value = 1398.0834
$ws.Range("J132").Value = 2506.625
$ws.Range("K132").Value = 4194.2502
$ws.Range("L132").Value = 7519.875
$ws.Range("M132").Value = -1664.2502
$ws.Range("N132").Value = -12579.875
$ws.Range("H136").Value = 3740.9524
$ws.Range("I136").Value = 1684.1111
$ws.Range("J136").Value = 5283.5835
$ws.Range("K136").Value = 5052.3333
$ws.Range("L136").Value = 15850.7505
$ws.Range("M136").Value = -2502.3333
$ws.Range("N136").Value = -20950.7505
$ws.Range("H137").Value = 60925
$ws.Range("J137").Value = 60925
$ws.Range("L137").Value = 60925
$ws.Range("N137").Value = -71125
$ws.Range("H139").Value = 74028.57000000001
$ws.Range("J139").Value = 74028.57000000001
$ws.Range("L139").Value = 74028.57000000001
$ws.Range("N139").Value = -84308.57000000001
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = -51070
